$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Updated Students"

# Headers
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "group"
$ws.Range("D1").Value = "Scholarship"
$ws.Range("E1").Value = "gpa"
$ws.Range("F1").Value = "faculty"
$ws.Range("G1").Value = "newScholarship"

# Row 2 - Ermek
$ws.Range("A2").Value = 21021
$ws.Range("B2").Value = "Ermek"
$ws.Range("C2").Value = "eng-2401"
$ws.Range("D2").Value = 41898
$ws.Range("E2").Value = 2.79
$ws.Range("F2").Value = "Engineering"
$ws.Range("G2").Value = 41898

# Row 3 - Khadisha
$ws.Range("A3").Value = 23023
$ws.Range("B3").Value = "Khadisha"
$ws.Range("C3").Value = "eco-2405"
$ws.Range("D3").Value = 41898
$ws.Range("E3").Value = 3.01
$ws.Range("F3").Value = "Economics"
$ws.Range("G3").Value = 41898

# Row 4 - Bekzhan
$ws.Range("A4").Value = 25025
$ws.Range("B4").Value = "Bekzhan"
$ws.Range("C4").Value = "phy-2415"
$ws.Range("D4").Value = 63000
$ws.Range("E4").Value = 2.83
$ws.Range("F4").Value = "Philosophy"
$ws.Range("G4").Value = 63000

# Row 5 - Arystan
$ws.Range("A5").Value = 27027
$ws.Range("B5").Value = "Arystan"
$ws.Range("C5").Value = "mrk-2424"
$ws.Range("D5").Value = 67000
$ws.Range("E5").Value = 3.14
$ws.Range("F5").Value = "Marketing"
$ws.Range("G5").Value = 67000

# Update selection to match target (G1 active cell)
$ws.Range("G1").Select() | Out-Null
